# ds4owd_precourse_survey.xlsx — fix CLI question formatting (KoboToolbox compatibility)
#
# The "survey" sheet had the CLI usage question's label text accidentally
# split across rows 36-39 (each row holding a fragment of the original
# multi-line label, plus the trailing rows getting misread as separate
# question rows). This script:
#   1. Rejoins the full label into C36 as a single-line string and
#      restores the "required" flag (D36 = "no") that had been lost.
#   2. Removes the three stray fragment rows (37-39) that held the
#      leftover pieces of the old multi-line text, shifting all
#      subsequent rows up so the sheet is contiguous again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# 1) Rebuild the CLI usage question label as a single line, and mark it
#    required again (matches the "no" used throughout the rest of the sheet).
$ws.Range("C36").Value = "Which of these best describes your current usage of the default command-line interface (CLI)? On Mac: The default CLI app is Terminal and the default shell is Zsh (you may also use Bash or other shells). On Windows: The default CLI app is Windows Terminal which can run Command Prompt PowerShell and Bash (via Windows Subsystem for Linux). How would you describe your experience?"
$ws.Range("D36").Value = "no"

# 2) Delete the three leftover fragment rows; everything below shifts up,
#    so the LLM/project/agreements questions close up the gap left by the
#    now-merged CLI question (dimension becomes A1:V56).
$ws.Rows("37:39").Delete()
